$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the old marketing-event value in G5
$ws.Range("G5").ClearContents()

# Insert the new marketing event and shift the rest down one row
$ws.Range("G6").Value = "post on hacker news"
$ws.Range("G7").Value = "start showing up in google search organically"
$ws.Range("G8").Value = "DA articles bringing traffic to site"
$ws.Range("G9").Value = "articles on dev sites about DA"
$ws.Range("G10").Value = "google adwords campaign"
$ws.Range("G11").Value = "bump up # facebbok ads"

# Update the active selection to reflect where the user left the cursor
$ws.Range("B4").Select()
